$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for columns D,I,J,K,L,M,N,O,P,Q for rows 279..375 (post-edit state).
# Rows 279-280 are unchanged; rows 281-373 take on the values that used to sit
# two rows above them; rows 374-375 are brand new rows appended at the bottom
# (carrying what used to be rows 372-373's data).
$data = @(
    @(44187, 'Primera', 870, 4500, 5000, 4747, '$/docena de matas', 'Región Metropolitana', 791, 6),
    @(44187, 'Segunda', 220, 4000, 4000, 4000, '$/docena de matas', 'Región Metropolitana', 667, 6),
    @(44187, 'Primera', 870, 4500, 5000, 4747, '$/docena de matas', 'Región Metropolitana', 791, 6),
    @(44187, 'Segunda', 220, 4000, 4000, 4000, '$/docena de matas', 'Región Metropolitana', 667, 6),
    @(44266, 'Primera', 570, 3000, 3500, 3149, '$/docena de matas', 'Región Metropolitana', 525, 6),
    @(44266, 'Segunda', 90, 2500, 2500, 2500, '$/docena de matas', 'Región Metropolitana', 417, 6),
    @(44277, 'Primera', 290, 3500, 4000, 3638, '$/docena de matas', 'Región Metropolitana', 606, 6),
    @(44277, 'Segunda', 60, 2500, 2500, 2500, '$/docena de matas', 'Región Metropolitana', 417, 6),
    @(44391, 'Primera', 70, 4500, 4500, 4500, '$/paquete', 'Región de Arica y Parinacota', 4500, 1),
    @(44202, 'Primera', 1340, 3500, 4500, 3914, '$/docena de matas', 'Región Metropolitana', 652, 6),
    @(44202, 'Segunda', 410, 3000, 3500, 3195, '$/docena de matas', 'Región Metropolitana', 532, 6),
    @(44371, 'Primera', 160, 3500, 4000, 3719, '$/paquete', 'Región de Arica y Parinacota', 3719, 1),
    @(44249, 'Primera', 180, 3500, 3500, 3500, '$/docena de matas', 'Provincia de Chacabuco', 583, 6),
    @(44249, 'Primera', 400, 3000, 3500, 3288, '$/docena de matas', 'Región Metropolitana', 548, 6),
    @(44249, 'Segunda', 120, 2500, 2500, 2500, '$/docena de matas', 'Región Metropolitana', 417, 6),
    @(44225, 'Primera', 1380, 3000, 3500, 3138, '$/docena de matas', 'Región Metropolitana', 523, 6),
    @(44225, 'Segunda', 370, 2000, 2500, 2338, '$/docena de matas', 'Región Metropolitana', 390, 6),
    @(44445, 'Primera', 80, 4500, 5000, 4719, '$/paquete', 'Región de Arica y Parinacota', 4719, 1),
    @(44259, 'Primera', 570, 3000, 4000, 3500, '$/docena de matas', 'Región Metropolitana', 583, 6),
    @(44259, 'Segunda', 230, 2500, 3000, 2761, '$/docena de matas', 'Región Metropolitana', 460, 6),
    @(44328, 'Primera', 60, 3300, 3500, 3433, '$/paquete', 'Región de Arica y Parinacota', 3433, 1),
    @(44561, 'Primera', 590, 3500, 4000, 3805, '$/docena de matas', 'Región Metropolitana', 634, 6),
    @(44561, 'Segunda', 260, 3000, 3000, 3000, '$/docena de matas', 'Región Metropolitana', 500, 6),
    @(44396, 'Primera', 130, 4500, 4500, 4500, '$/paquete', 'Región de Arica y Parinacota', 4500, 1),
    @(44494, 'Primera', 78, 4500, 5000, 4705, '$/paquete', 'Región de Arica y Parinacota', 4705, 1),
    @(44526, 'Primera', 250, 5000, 6000, 5400, '$/docena de matas', 'Región Metropolitana', 900, 6),
    @(44526, 'Segunda', 100, 4000, 4000, 4000, '$/docena de matas', 'Región Metropolitana', 667, 6),
    @(44250, 'Primera', 160, 3000, 3000, 3000, '$/docena de matas', 'Provincia de Chacabuco', 500, 6),
    @(44250, 'Primera', 1050, 3000, 3500, 3214, '$/docena de matas', 'Región Metropolitana', 536, 6),
    @(44250, 'Segunda', 80, 2000, 2000, 2000, '$/docena de matas', 'Provincia de Chacabuco', 333, 6),
    @(44250, 'Segunda', 570, 2000, 2500, 2395, '$/docena de matas', 'Región Metropolitana', 399, 6),
    @(44285, 'Primera', 110, 3000, 3000, 3000, '$/docena de matas', 'Región Metropolitana', 500, 6),
    @(44285, 'Segunda', 150, 2500, 2500, 2500, '$/docena de matas', 'Región Metropolitana', 417, 6),
    @(44264, 'Primera', 90, 3000, 3000, 3000, '$/docena de matas', 'Provincia de Chacabuco', 500, 6),
    @(44264, 'Primera', 150, 3000, 3000, 3000, '$/docena de matas', 'Región Metropolitana', 500, 6),
    @(44264, 'Segunda', 50, 2500, 2500, 2500, '$/docena de matas', 'Provincia de Chacabuco', 417, 6),
    @(44264, 'Segunda', 80, 2500, 2500, 2500, '$/docena de matas', 'Región Metropolitana', 417, 6),
    @(44533, 'Primera', 620, 4500, 6000, 5113, '$/docena de matas', 'Región Metropolitana', 852, 6),
    @(44533, 'Segunda', 180, 4000, 4500, 4222, '$/docena de matas', 'Región Metropolitana', 704, 6),
    @(44489, 'Primera', 150, 5000, 6000, 5467, '$/paquete', 'Región de Arica y Parinacota', 5467, 1),
    @(44354, 'Primera', 140, 3000, 3500, 3286, '$/paquete', 'Región de Arica y Parinacota', 3286, 1),
    @(44221, 'Primera', 220, 3000, 3000, 3000, '$/docena de matas', 'Región Metropolitana', 500, 6),
    @(44221, 'Segunda', 140, 2500, 2500, 2500, '$/docena de matas', 'Región Metropolitana', 417, 6),
    @(44523, 'Primera', 550, 4500, 5000, 4791, '$/docena de matas', 'Región Metropolitana', 798, 6),
    @(44523, 'Segunda', 250, 4000, 4000, 4000, '$/docena de matas', 'Región Metropolitana', 667, 6),
    @(44399, 'Primera', 130, 4500, 4500, 4500, '$/paquete', 'Región de Arica y Parinacota', 4500, 1),
    @(44382, 'Primera', 70, 4000, 4500, 4286, '$/paquete', 'Región de Arica y Parinacota', 4286, 1),
    @(44441, 'Primera', 70, 4500, 5000, 4836, '$/paquete', 'Región de Arica y Parinacota', 4836, 1),
    @(44167, 'Primera', 520, 4500, 5000, 4885, '$/docena de matas', 'Región Metropolitana', 814, 6),
    @(44167, 'Segunda', 130, 4000, 4000, 4000, '$/docena de matas', 'Región Metropolitana', 667, 6),
    @(44335, 'Primera', 100, 3000, 3500, 3300, '$/paquete', 'Región de Arica y Parinacota', 3300, 1),
    @(44316, 'Primera', 130, 5000, 5000, 5000, '$/docena de matas', 'Región Metropolitana', 833, 6),
    @(44316, 'Segunda', 70, 4000, 4000, 4000, '$/docena de matas', 'Región Metropolitana', 667, 6),
    @(44475, 'Primera', 77, 4000, 5000, 4416, '$/paquete', 'Región de Arica y Parinacota', 4416, 1),
    @(44315, 'Primera', 220, 5000, 5000, 5000, '$/docena de matas', 'Región Metropolitana', 833, 6),
    @(44315, 'Segunda', 140, 4000, 4000, 4000, '$/docena de matas', 'Región Metropolitana', 667, 6),
    @(44186, 'Primera', 820, 4000, 5000, 4613, '$/docena de matas', 'Región Metropolitana', 769, 6),
    @(44186, 'Segunda', 230, 3000, 4000, 3696, '$/docena de matas', 'Región Metropolitana', 616, 6),
    @(44176, 'Primera', 1280, 4500, 5000, 4824, '$/docena de matas', 'Región Metropolitana', 804, 6),
    @(44176, 'Segunda', 80, 4000, 4000, 4000, '$/docena de matas', 'Región Metropolitana', 667, 6),
    @(44496, 'Primera', 100, 4500, 5000, 4800, '$/paquete', 'Región de Arica y Parinacota', 4800, 1),
    @(44278, 'Primera', 590, 3500, 4000, 3720, '$/docena de matas', 'Región Metropolitana', 620, 6),
    @(44278, 'Segunda', 80, 2500, 2500, 2500, '$/docena de matas', 'Región Metropolitana', 417, 6),
    @(44438, 'Primera', 200, 5000, 5500, 5300, '$/paquete', 'Región de Arica y Parinacota', 5300, 1),
    @(44312, 'Primera', 70, 5000, 5000, 5000, '$/docena de matas', 'Región Metropolitana', 833, 6),
    @(44312, 'Segunda', 150, 4000, 4000, 4000, '$/docena de matas', 'Región Metropolitana', 667, 6),
    @(44300, 'Segunda', 130, 4000, 4000, 4000, '$/docena de matas', 'Región Metropolitana', 667, 6),
    @(44314, 'Primera', 140, 5000, 5000, 5000, '$/docena de matas', 'Región Metropolitana', 833, 6),
    @(44314, 'Segunda', 90, 4000, 4000, 4000, '$/docena de matas', 'Región Metropolitana', 667, 6),
    @(44260, 'Primera', 620, 3500, 4500, 3960, '$/docena de matas', 'Región Metropolitana', 660, 6),
    @(44260, 'Segunda', 230, 3000, 3000, 3000, '$/docena de matas', 'Región Metropolitana', 500, 6),
    @(44585, 'Primera', 520, 3000, 3500, 3240, '$/docena de matas', 'Región Metropolitana', 540, 6),
    @(44585, 'Segunda', 180, 2000, 2500, 2222, '$/docena de matas', 'Región Metropolitana', 370, 6),
    @(44560, 'Primera', 1130, 4000, 4500, 4155, '$/docena de matas', 'Región Metropolitana', 692, 6),
    @(44560, 'Segunda', 410, 3000, 3500, 3305, '$/docena de matas', 'Región Metropolitana', 551, 6),
    @(44272, 'Primera', 420, 3000, 3500, 3179, '$/docena de matas', 'Región Metropolitana', 530, 6),
    @(44272, 'Segunda', 80, 2500, 2500, 2500, '$/docena de matas', 'Región Metropolitana', 417, 6),
    @(44385, 'Primera', 60, 4500, 5000, 4667, '$/paquete', 'Región de Arica y Parinacota', 4667, 1),
    @(44162, 'Primera', 410, 4500, 5000, 4902, '$/docena de matas', 'Región Metropolitana', 817, 6),
    @(44162, 'Segunda', 210, 4000, 4000, 4000, '$/docena de matas', 'Región Metropolitana', 667, 6),
    @(44529, 'Primera', 130, 6000, 6000, 6000, '$/docena de matas', 'Región Metropolitana', 1000, 6),
    @(44529, 'Segunda', 60, 5000, 5000, 5000, '$/docena de matas', 'Región Metropolitana', 833, 6),
    @(44323, 'Primera', 130, 5000, 5000, 5000, '$/docena de matas', 'Región Metropolitana', 833, 6),
    @(44323, 'Segunda', 100, 4000, 4000, 4000, '$/docena de matas', 'Región Metropolitana', 667, 6),
    @(44306, 'Primera', 160, 5000, 5000, 5000, '$/docena de matas', 'Región Metropolitana', 833, 6),
    @(44189, 'Primera', 1170, 4000, 5000, 4380, '$/docena de matas', 'Región Metropolitana', 730, 6),
    @(44189, 'Segunda', 410, 3500, 4000, 3695, '$/docena de matas', 'Región Metropolitana', 616, 6),
    @(44321, 'Primera', 130, 5000, 5000, 5000, '$/docena de matas', 'Región Metropolitana', 833, 6),
    @(44321, 'Segunda', 60, 4000, 4000, 4000, '$/docena de matas', 'Región Metropolitana', 667, 6),
    @(44302, 'Primera', 190, 6000, 6000, 6000, '$/docena de matas', 'Región Metropolitana', 1000, 6),
    @(44209, 'Primera', 690, 3500, 4000, 3710, '$/docena de matas', 'Región Metropolitana', 618, 6),
    @(44209, 'Segunda', 310, 3000, 3000, 3000, '$/docena de matas', 'Región Metropolitana', 500, 6),
    @(44274, 'Primera', 720, 3000, 25000, 11139, '$/docena de matas', 'Región Metropolitana', 1856, 6),
    @(44274, 'Segunda', 90, 3000, 3000, 3000, '$/docena de matas', 'Región Metropolitana', 500, 6),
    @(44554, 'Primera', 570, 4000, 5500, 4833, '$/docena de matas', 'Región Metropolitana', 806, 6),
    @(44554, 'Segunda', 210, 3500, 4500, 3976, '$/docena de matas', 'Región Metropolitana', 663, 6),
    @(44392, 'Primera', 70, 4000, 4000, 4000, '$/paquete', 'Región de Arica y Parinacota', 4000, 1)
)

$startRow = 279
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 4).Value = $row[0]   # D Fecha
    $ws.Cells.Item($r, 9).Value = $row[1]   # I Calidad
    $ws.Cells.Item($r, 10).Value = $row[2]  # J Volumen
    $ws.Cells.Item($r, 11).Value = $row[3]  # K Precio minimo
    $ws.Cells.Item($r, 12).Value = $row[4]  # L Precio maximo
    $ws.Cells.Item($r, 13).Value = $row[5]  # M Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $row[6]  # N Unidad de comercializacion
    $ws.Cells.Item($r, 15).Value = $row[7]  # O Origen
    $ws.Cells.Item($r, 16).Value = $row[8]  # P Precio $/Kg
    $ws.Cells.Item($r, 17).Value = $row[9]  # Q Kg o Unidades
}

# The two brand-new rows (374, 375) need the constant columns too -- these are
# identical across the whole table (A, B, C, E, F, G, H, R).
$ws.Cells.Item(374, 1).Value = 6
$ws.Cells.Item(374, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(374, 3).Value = "Metropolitana"
$ws.Cells.Item(374, 5).Value = 13
$ws.Cells.Item(374, 6).Value = 100112052
$ws.Cells.Item(374, 7).Value = "Albahaca"
$ws.Cells.Item(374, 8).Value = "Sin especificar"
$ws.Cells.Item(374, 18).Value = "Hortaliza"

$ws.Cells.Item(375, 1).Value = 6
$ws.Cells.Item(375, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(375, 3).Value = "Metropolitana"
$ws.Cells.Item(375, 5).Value = 13
$ws.Cells.Item(375, 6).Value = 100112052
$ws.Cells.Item(375, 7).Value = "Albahaca"
$ws.Cells.Item(375, 8).Value = "Sin especificar"
$ws.Cells.Item(375, 18).Value = "Hortaliza"

# Copy the date style (style index 2, applied to column D) down onto the new rows.
$ws.Range("D373").Copy()
$ws.Range("D374:D375").PasteSpecial(-4122)
